$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.701.67'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.72%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.167.04'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.96%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '616.79'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.29%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '148.46'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.74%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.165.64'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.97%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.530'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.68%  '
$ws.Range('E10').Value = '  -0.69%  '
$ws.Range('E11').Value = '  -2.06%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.474'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.86%  '
$ws.Range('E13').Value = '  +0.06%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.89'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.95%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.692.47'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.18%  '
$ws.Range('E16').Value = '  +2.85%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.726.40'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.65%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.168.43'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.38%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.94'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.61%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '480.09'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.73%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.76'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.725'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.63%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.98'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.58%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.80'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.87%  '
$ws.Range('E25').Value = '  -0.30%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('E27').Value = '  -2.81%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.60'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.07%  '
$ws.Range('E29').Value = '  -1.27%  '
$ws.Range('E30').Value = '  -6.01%  '
$ws.Range('E31').Value = '  -7.65%  '
$ws.Range('E32').Value = '  +0.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.71'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.13%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '26.58'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.00%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.14'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.18%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0₃0779'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.96%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.01'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.71%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '53.07'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.96%  '
$ws.Range('B39').Value = 'dogwifhat'
$ws.Range('C39').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.18'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.19%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '461.74'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.51%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0400'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.84%  '
$ws.Range('E42').Value = '  -3.96%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.43'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.22%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.847.55'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.48%  '
$ws.Range('E45').Value = '  -3.41%  '
$ws.Range('E46').Value = '  -2.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.44'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.58%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '26.69'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.08%  '
$ws.Range('E50').Value = '  -1.50%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '120.38'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.14%  '
